# Bug Metrics v1.xlsx - "Debugged and Updated Bug Metric" edit
#
# Displays the dates of vital-signs charts as plain text/string on the
# "Iteration 5" sheet so the chart renders consistently across browsers
# (IE / Chrome / Openshift) instead of converting a native Date, and
# records who made the change + when it was solved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iteration 5")

# Rows 9-13 ("Action Taken by Developers" / "Solved by" / "Date Solved")
# were previously blank - fill them in.
$note   = "Printed the dates of vital signs as String on chart instead of Date to avoid conversion error at different platforms"
$solver = "Wei Yi"
$dateSolved = 42002   # 29-Dec-2014 (serial date, matches the "Date Solved" column format used elsewhere on the sheet)

foreach ($r in 9..13) {
    $ws.Cells.Item($r, 10).Value = $note        # column J - Action Taken by Developers
    $ws.Cells.Item($r, 11).Value = $solver       # column K - Solved by
    $ws.Cells.Item($r, 12).Value = $dateSolved   # column L - Date Solved
}

# Rows 11-13 now need extra height to fit the wrapped note text (rows 9
# and 10 already had a custom height tall enough).
$ws.Cells.Item(11, 10).EntireRow.RowHeight = 54.75
$ws.Cells.Item(12, 10).EntireRow.RowHeight = 54.75
$ws.Cells.Item(13, 10).EntireRow.RowHeight = 54.75

# Move the active selection from H11 to J11 (the column that was just edited).
$ws.Activate() | Out-Null
$ws.Range("J11").Select() | Out-Null
